$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for "Macroferia Regional de Talca -
# Pepino ensalada". Insert a fresh row at 123, which pushes the existing
# rows 123-204 down to 124-205, then populate the new row with its data.
$ws.Rows.Item(123).Insert()

$ws.Range("A123").Value = 5
$ws.Range("B123").Value = "Macroferia Regional de Talca"
$ws.Range("C123").Value = "Maule"
$ws.Range("D123").Value = 44438
$ws.Range("E123").Value = 7
$ws.Range("F123").Value = 100112043
$ws.Range("G123").Value = "Pepino ensalada"
$ws.Range("H123").Value = "Sin especificar"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 400
$ws.Range("K123").Value = 13000
$ws.Range("L123").Value = 13000
$ws.Range("M123").Value = 13000
$ws.Range("N123").Value = "`$/caja 60 unidades"
$ws.Range("O123").Value = "Región del Maule"
$ws.Range("P123").Value = 217
$ws.Range("Q123").Value = 60
$ws.Range("R123").Value = "Hortaliza"
